$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPLKAKT064-002")

# Regresi Tanggal: update the Periode Harian verification date in R2
$ws.Range("R2").Value = 20240907

# Move/save the active selection to R3 (matches the saved cursor position)
$ws.Range("R3").Select() | Out-Null
